$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 4 de Septiembre de 2020 a las 01:00"

# Row 4
$ws.Range("B4").Value = 6331042
$ws.Range("C4").Value = 40305
$ws.Range("D4").Value = 3571036
$ws.Range("E4").Value = 2569048
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 994
$ws.Range("H4").Value = 190958

# Row 5
$ws.Range("B5").Value = 4041638
$ws.Range("C5").Value = 40216
$ws.Range("D5").Value = 3247610
$ws.Range("E5").Value = 669377
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 752
$ws.Range("H5").Value = 124651

# Row 9
$ws.Range("B9").Value = 641574
$ws.Range("C9").Value = 8235
$ws.Range("D9").Value = 489151
$ws.Range("E9").Value = 131805
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 270
$ws.Range("H9").Value = 20618

# Row 13
$ws.Range("B13").Value = 451198
$ws.Range("C13").Value = 12026
$ws.Range("D13").Value = 322461
$ws.Range("E13").Value = 119376
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 243
$ws.Range("H13").Value = 9361

# Row 27
$ws.Range("B27").Value = 130493
$ws.Range("C27").Value = 570
$ws.Range("D27").Value = 115444
$ws.Range("E27").Value = 5908
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 9141

# Row 34
$ws.Range("B34").Value = 99425
$ws.Range("C34").Value = 145
$ws.Range("D34").Value = 75415
$ws.Range("E34").Value = 18531
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 18
$ws.Range("H34").Value = 5479

# Row 35
$ws.Range("B35").Value = 96629
$ws.Range("C35").Value = 1002
$ws.Range("D35").Value = 70871
$ws.Range("E35").Value = 23957
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 1765

# Row 47 -> Japon
$ws.Range("A47").Value = "Japon"
$ws.Range("B47").Value = 69599
$ws.Range("C47").Value = 598
$ws.Range("D47").Value = 59524
$ws.Range("E47").Value = 8756
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 12
$ws.Range("H47").Value = 1319

# Row 48 -> Polonia
$ws.Range("A48").Value = "Polonia"
$ws.Range("B48").Value = 69129
$ws.Range("C48").Value = 612
$ws.Range("D48").Value = 48593
$ws.Range("E48").Value = 18444
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 14
$ws.Range("H48").Value = 2092

# Row 54
$ws.Range("B54").Value = 54587
$ws.Range("C54").Value = 124
$ws.Range("D54").Value = 42622
$ws.Range("E54").Value = 10917
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = 21
$ws.Range("H54").Value = 1048

# Row 55
$ws.Range("B55").Value = 53433
$ws.Range("C55").Value = 626
$ws.Range("D55").Value = 50013
$ws.Range("E55").Value = 3230
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 0
$ws.Range("H55").Value = 190

# Row 78
$ws.Range("B78").Value = 19604
$ws.Range("C78").Value = 144
$ws.Range("D78").Value = 18448
$ws.Range("E78").Value = 741
$ws.Range("F78").Value = 0
$ws.Range("G78").Value = 0
$ws.Range("H78").Value = 415

# Row 90 -> Noruega
$ws.Range("A90").Value = "Noruega"
$ws.Range("B90").Value = 11120
$ws.Range("C90").Value = 86
$ws.Range("D90").Value = 9348
$ws.Range("E90").Value = 1508
$ws.Range("F90").Value = 0
$ws.Range("G90").Value = 0
$ws.Range("H90").Value = 264

# Row 91 -> Croacia
$ws.Range("A91").Value = "Croacia"
$ws.Range("B91").Value = 11094
$ws.Range("C91").Value = 369
$ws.Range("D91").Value = 8266
$ws.Range("E91").Value = 2634
$ws.Range("F91").Value = 0
$ws.Range("G91").Value = 3
$ws.Range("H91").Value = 194

# Row 106
$ws.Range("B106").Value = 6811
$ws.Range("C106").Value = 66
$ws.Range("D106").Value = 0
$ws.Range("E106").Value = 0
$ws.Range("F106").Value = 0
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 124

# Row 110
$ws.Range("B110").Value = 5165
$ws.Range("C110").Value = 146
$ws.Range("D110").Value = 4164
$ws.Range("E110").Value = 897
$ws.Range("F110").Value = 0
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 104

# Row 145 -> Trinidad yTobago
$ws.Range("A145").Value = "Trinidad yTobago"
$ws.Range("B145").Value = 1984
$ws.Range("C145").Value = 64
$ws.Range("D145").Value = 700
$ws.Range("E145").Value = 1255
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 1
$ws.Range("H145").Value = 29

# Row 146 -> Yemen
$ws.Range("A146").Value = "Yemen"
$ws.Range("B146").Value = 1979
$ws.Range("C146").Value = 3
$ws.Range("D146").Value = 1180
$ws.Range("E146").Value = 228
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 571

# Row 147 -> Malta
$ws.Range("A147").Value = "Malta"
$ws.Range("B147").Value = 1965
$ws.Range("C147").Value = 34
$ws.Range("D147").Value = 1528
$ws.Range("E147").Value = 424
$ws.Range("F147").Value = 0
$ws.Range("G147").Value = 0
$ws.Range("H147").Value = 13

# Row 151
$ws.Range("B151").Value = 1636
$ws.Range("C151").Value = 10
$ws.Range("D151").Value = 1437
$ws.Range("E151").Value = 155
$ws.Range("F151").Value = 0
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 44

# Row 156
$ws.Range("B156").Value = 1401
$ws.Range("C156").Value = 19
$ws.Range("D156").Value = 805
$ws.Range("E156").Value = 552
$ws.Range("F156").Value = 0
$ws.Range("G156").Value = 3
$ws.Range("H156").Value = 44

# Row 166
$ws.Range("B166").Value = 1018
$ws.Range("C166").Value = 1
$ws.Range("D166").Value = 908
$ws.Range("E166").Value = 33
$ws.Range("F166").Value = 0
$ws.Range("G166").Value = 0
$ws.Range("H166").Value = 77

# Row 167
$ws.Range("B167").Value = 897
$ws.Range("C167").Value = 1
$ws.Range("D167").Value = 857
$ws.Range("E167").Value = 25
$ws.Range("F167").Value = 0
$ws.Range("G167").Value = 0
$ws.Range("H167").Value = 15

# Row 178
$ws.Range("B178").Value = 448
$ws.Range("C178").Value = 21
$ws.Range("D178").Value = 410
$ws.Range("E178").Value = 31
$ws.Range("F178").Value = 0
$ws.Range("G178").Value = 0
$ws.Range("H178").Value = 7

# Row 214 -> Montserrat
$ws.Range("A214").Value = "Montserrat"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 12
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 1

# Row 215 -> Islas Malvinas
$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("B215").Value = 13
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 13
$ws.Range("E215").Value = 0
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 0
